# Natmi following Dr Hou advice
# Rewrites Sheet1 of the Tslp-Il7r LR-pair table with the updated
# natmi output: 4 sending clusters (ECs, FAPs, M2, sCs) each paired with
# 2 target clusters (ECs, M2) for the Tslp -> Il7r interaction (rows 2-9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns (A-D) are written column-by-column so that new shared
# strings are interned in a stable, predictable order.
# Column A: Sending cluster
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("A4").Value = "FAPs"
$ws.Range("A5").Value = "FAPs"
$ws.Range("A6").Value = "M2"
$ws.Range("A7").Value = "M2"
$ws.Range("A8").Value = "sCs"
$ws.Range("A9").Value = "sCs"

# Column B: Ligand symbol
$ws.Range("B2").Value = "Tslp"
$ws.Range("B3").Value = "Tslp"
$ws.Range("B4").Value = "Tslp"
$ws.Range("B5").Value = "Tslp"
$ws.Range("B6").Value = "Tslp"
$ws.Range("B7").Value = "Tslp"
$ws.Range("B8").Value = "Tslp"
$ws.Range("B9").Value = "Tslp"

# Column C: Receptor symbol
$ws.Range("C2").Value = "Il7r"
$ws.Range("C3").Value = "Il7r"
$ws.Range("C4").Value = "Il7r"
$ws.Range("C5").Value = "Il7r"
$ws.Range("C6").Value = "Il7r"
$ws.Range("C7").Value = "Il7r"
$ws.Range("C8").Value = "Il7r"
$ws.Range("C9").Value = "Il7r"

# Column D: Target cluster
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "M2"
$ws.Range("D4").Value = "ECs"
$ws.Range("D5").Value = "M2"
$ws.Range("D6").Value = "ECs"
$ws.Range("D7").Value = "M2"
$ws.Range("D8").Value = "ECs"
$ws.Range("D9").Value = "M2"

# Numeric columns (E-T): natmi-computed statistics for each row.
# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.158268666666667
$ws.Range("H2").Value = 3.474806
$ws.Range("I2").Value = 0.3523202827966646
$ws.Range("J2").Value = 0.3523202827966647
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 3.970254666666667
$ws.Range("N2").Value = 11.910764
$ws.Range("O2").Value = 0.1277192879665705
$ws.Range("P2").Value = 0.1277192879665705
$ws.Range("Q2").Value = 4.598621579087111
$ws.Range("R2").Value = 41.38759421178401
$ws.Range("S2").Value = 0.04499809565497077
$ws.Range("T2").Value = 0.04499809565497078

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.158268666666667
$ws.Range("H3").Value = 3.474806
$ws.Range("I3").Value = 0.3523202827966646
$ws.Range("J3").Value = 0.3523202827966647
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 27.115533
$ws.Range("N3").Value = 81.34659900000001
$ws.Range("O3").Value = 0.8722807120334295
$ws.Range("P3").Value = 0.8722807120334295
$ws.Range("Q3").Value = 31.407072253866
$ws.Range("R3").Value = 282.6636502847941
$ws.Range("S3").Value = 0.3073221871416938
$ws.Range("T3").Value = 0.3073221871416939

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8431363333333334
$ws.Range("H4").Value = 2.529409
$ws.Range("I4").Value = 0.2564638412010423
$ws.Range("J4").Value = 0.2564638412010423
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 3.970254666666667
$ws.Range("N4").Value = 11.910764
$ws.Range("O4").Value = 0.1277192879665705
$ws.Range("P4").Value = 0.1277192879665705
$ws.Range("Q4").Value = 3.347465962052889
$ws.Range("R4").Value = 30.127193658476
$ws.Range("S4").Value = 0.03275537918736873
$ws.Range("T4").Value = 0.03275537918736873

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.8431363333333334
$ws.Range("H5").Value = 2.529409
$ws.Range("I5").Value = 0.2564638412010423
$ws.Range("J5").Value = 0.2564638412010423
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 27.115533
$ws.Range("N5").Value = 81.34659900000001
$ws.Range("O5").Value = 0.8722807120334295
$ws.Range("P5").Value = 0.8722807120334295
$ws.Range("Q5").Value = 22.862091069999
$ws.Range("R5").Value = 205.7588196299911
$ws.Range("S5").Value = 0.2237084620136735
$ws.Range("T5").Value = 0.2237084620136735

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.749088
$ws.Range("H6").Value = 2.247264
$ws.Range("I6").Value = 0.2278563718373814
$ws.Range("J6").Value = 0.2278563718373814
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 3.970254666666667
$ws.Range("N6").Value = 11.910764
$ws.Range("O6").Value = 0.1277192879665705
$ws.Range("P6").Value = 0.1277192879665705
$ws.Range("Q6").Value = 2.974070127744
$ws.Range("R6").Value = 26.766631149696
$ws.Range("S6").Value = 0.02910165356971648
$ws.Range("T6").Value = 0.02910165356971648

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.749088
$ws.Range("H7").Value = 2.247264
$ws.Range("I7").Value = 0.2278563718373814
$ws.Range("J7").Value = 0.2278563718373814
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 27.115533
$ws.Range("N7").Value = 81.34659900000001
$ws.Range("O7").Value = 0.8722807120334295
$ws.Range("P7").Value = 0.8722807120334295
$ws.Range("Q7").Value = 20.311920383904
$ws.Range("R7").Value = 182.807283455136
$ws.Range("S7").Value = 0.1987547182676649
$ws.Range("T7").Value = 0.1987547182676649

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.5370516666666667
$ws.Range("H8").Value = 1.611155
$ws.Range("I8").Value = 0.1633595041649117
$ws.Range("J8").Value = 0.1633595041649117
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 3.970254666666667
$ws.Range("N8").Value = 11.910764
$ws.Range("O8").Value = 0.1277192879665705
$ws.Range("P8").Value = 0.1277192879665705
$ws.Range("Q8").Value = 2.132231885824445
$ws.Range("R8").Value = 19.19008697242
$ws.Range("S8").Value = 0.02086415955451454
$ws.Range("T8").Value = 0.02086415955451454

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.5370516666666667
$ws.Range("H9").Value = 1.611155
$ws.Range("I9").Value = 0.1633595041649117
$ws.Range("J9").Value = 0.1633595041649117
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 27.115533
$ws.Range("N9").Value = 81.34659900000001
$ws.Range("O9").Value = 0.8722807120334295
$ws.Range("P9").Value = 0.8722807120334295
$ws.Range("Q9").Value = 14.562442190205
$ws.Range("R9").Value = 131.061979711845
$ws.Range("S9").Value = 0.1424953446103972
$ws.Range("T9").Value = 0.1424953446103972
